# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value, per the diff
$updates = @{
    "F5"  = -1
    "E9"  = -1
    "H9"  = 0
    "I9"  = 8
    "F11" = -2
    "F13" = -1
    "F21" = -1
    "F23" = 3
    "F24" = -2
    "F26" = 4
    "F37" = -5
    "F44" = -3
    "F49" = -5
    "F53" = -1
    "F54" = -3
    "F58" = -2
    "F60" = -2
    "F69" = 1
    "F70" = 0
    "F73" = -5
    "F75" = -10
    "F76" = -6
    "F77" = 7
    "F78" = -2
    "F80" = -4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
